$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# 1. Fix "neévre" typo -> "névre" and de-duplicate the
#    "hallgató bounty hunter" phrase that was previously split across a
#    second, redundant set of runs.
Replace-Text "Wild West Wind neévre hallgató bounty hunter szerű" `
             "Wild West Wind névre hallgató bounty hunter szerű"

# 2. Drop the stray trailing space run at the end of the paragraph.
Replace-Text "kezdőhelyzete ez: " "kezdőhelyzete ez:"

# 3. Add commas around the relative clauses.
Replace-Text "Jake Colton egy cowboy aki egyedül járja a környéket a lovával akit Rustlernek hívnak" `
             "Jake Colton egy cowboy, aki egyedül járja a környéket a lovával, akit Rustlernek hívnak"

# 4. Fix "Belékötnek" -> "Belekötnek".
Replace-Text "Belékötnek és megverik" "Belekötnek és megverik"

# 5. Rewrite the "fő cél" sentence.
Replace-Text "A játék során a fő cél a Desert Vultures banda tagjainak egyenkénti levadászása lesz a játékos célja és a bandavezér Mad dog McCoy megölése mint aféle fő boss." `
             "A játék során a cél a banda tagjainak egyenkénti levadászása lesz, a játékos fő célja a bandavezér, Mad dog McCoy megölése lesz, mint aféle fő boss."

# 6. Add commas in the "kalandozhat a játéktérképen is" sentence.
Replace-Text "kalandozhat a játéktérképen is ami 5 különböző részből vagyis biom-ból áll." `
             "kalandozhat a játéktérképen is, ami 5 különböző részből, vagyis biom-ból áll."

# 7. Add comma before "hogy" in the "jelenti" sentence.
Replace-Text "banda tagjait, ami azt jelenti hogy a játékosnak" `
             "banda tagjait, ami azt jelenti, hogy a játékosnak"

# 8. "A játékérték:" -> "A játéktérkép:"
Replace-Text "A játékérték:" "A játéktérkép:"
